$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the three rows get re-ordered by File Name, and the
# bce2a8a9 row (now row 4) gets its status flipped to "Ready for handoff"
# with fresh timestamps.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$wsOverview.Range("B2").Value = "e2e\ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$wsOverview.Range("G2").Value = "2016-08-30 19:17:37"

$wsOverview.Range("A3").Value = "fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$wsOverview.Range("B3").Value = "e2e\fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$wsOverview.Range("G3").Value = "2016-08-30 19:17:37"

$wsOverview.Range("A4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsOverview.Range("B4").Value = "e2e\bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-30 19:20:29"

$ovLinks = @($wsOverview.Hyperlinks)
$ovLinks[0].TextToDisplay = "e2e\ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$ovLinks[1].TextToDisplay = "e2e\fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$ovLinks[2].TextToDisplay = "e2e\bce2a8a9-9065-499c-bf50-4df3092b2d37.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row re-ordering; bce2a8a9 row (row 4) becomes
# "Ready for handoff" with its own handoff file/date and a new error detail
# describing a stale handback version.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$wsZh.Range("G2").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-30 19:17:31"
$wsZh.Range("I2").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$wsZh.Range("J2").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 19:17:54"

$wsZh.Range("A3").Value = "fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-30 19:17:31"
$wsZh.Range("I3").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$wsZh.Range("J3").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-30 19:17:54"

$wsZh.Range("A4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("G4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.6c72555e93cfaabf718348890bba91a9a180c88d.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-30 19:20:23"
$wsZh.Range("I4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsZh.Range("J4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.6c72555e93cfaabf718348890bba91a9a180c88d.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-08-30 19:19:39"
$wsZh.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07b0780518c70fb4dd3dc2d24b9f368ab50170e5/e2e/bce2a8a9-9065-499c-bf50-4df3092b2d37.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f553099104000f24d9d90208dd9b5fdcafb69617/e2e/bce2a8a9-9065-499c-bf50-4df3092b2d37.md."

$wsZh.Columns.Item(16).ColumnWidth = 39.17

$zhLinks = @($wsZh.Hyperlinks)
$zhLinks[0].TextToDisplay = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$zhLinks[1].TextToDisplay = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$zhLinks[2].TextToDisplay = "fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$zhLinks[3].TextToDisplay = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$zhLinks[4].TextToDisplay = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$zhLinks[5].TextToDisplay = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"

# ---------------------------------------------------------------------------
# Sheet "de-de": mirrors the "zh-cn" sheet treatment but keeps its own
# de-de file names/dates.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$wsDe.Range("H2").Value = "2016-08-30 19:17:37"
$wsDe.Range("I2").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$wsDe.Range("K2").Value = "2016-08-30 19:18:03"

$wsDe.Range("A3").Value = "fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("H3").Value = "2016-08-30 19:17:37"
$wsDe.Range("I3").Value = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$wsDe.Range("K3").Value = "2016-08-30 19:18:03"

$wsDe.Range("A4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("G4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.6c72555e93cfaabf718348890bba91a9a180c88d.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-30 19:20:29"
$wsDe.Range("I4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$wsDe.Range("J4").Value = "bce2a8a9-9065-499c-bf50-4df3092b2d37.6c72555e93cfaabf718348890bba91a9a180c88d.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-30 19:19:47"
$wsDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07b0780518c70fb4dd3dc2d24b9f368ab50170e5/e2e/bce2a8a9-9065-499c-bf50-4df3092b2d37.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f553099104000f24d9d90208dd9b5fdcafb69617/e2e/bce2a8a9-9065-499c-bf50-4df3092b2d37.md."

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$deLinks = @($wsDe.Hyperlinks)
$deLinks[0].TextToDisplay = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"
$deLinks[1].TextToDisplay = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$deLinks[2].TextToDisplay = "fffffff12a1cc5-c91b-47f5-9129-2deb057c1712.md"
$deLinks[3].TextToDisplay = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"
$deLinks[4].TextToDisplay = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"
$deLinks[5].TextToDisplay = "bce2a8a9-9065-499c-bf50-4df3092b2d37.md"

Write-Output "Report regenerated for handoff"
